$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "\\192.168.29.200\f\service\0.样品管理部\01 蛋白库相关\06 蛋白编号\00 理化质检-P90000之后在这里查理化质检结果\SEC"
$ws.Range("B2").Value = "【SEC】WKL230904-1 ZJ004 HLX1005 HNKJ001.pptx"
$ws.Range("C2").Value = "File is not a zip file"

$ws.Range("A3").Value = "\\192.168.29.200\f\service\0.样品管理部\01 蛋白库相关\06 蛋白编号\00 理化质检-P90000之后在这里查理化质检结果\SEC"
$ws.Range("B3").Value = "【SEC】WKL230904-1 ZJ004 HLX1005 HNKJ001.pdf"
$ws.Range("C3").Value = "NoRelatedPPT"
